$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - update "想去人数" (F column) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 7250
$ws1.Range("F4").Value = 5428
$ws1.Range("F8").Value = 42
$ws1.Range("F10").Value = 82
$ws1.Range("F12").Value = 201
$ws1.Range("F13").Value = 18
$ws1.Range("F14").Value = 645
$ws1.Range("F15").Value = 253
$ws1.Range("F18").Value = 20

# Sheet "全部类型" (sheet4) - update "想去人数" (F column) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 7251
$ws4.Range("F4").Value = 5428
$ws4.Range("F8").Value = 42
$ws4.Range("F10").Value = 82
$ws4.Range("F12").Value = 201
$ws4.Range("F13").Value = 18
$ws4.Range("F14").Value = 645
$ws4.Range("F15").Value = 253
$ws4.Range("F18").Value = 20
